$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 21.06924211128599
$ws.Range("C2").Value = 13.58868446670183
$ws.Range("D2").Value = 6.372331988624137
$ws.Range("F2").Value = 47.67677978485337
$ws.Range("G2").Value = 3.718125903691642
$ws.Range("J2").Value = 10.46194160280927
$ws.Range("M2").Value = 20.82179846458828
$ws.Range("N2").Value = 20.71968954315292

$ws.Range("B3").Value = 20.63667133388191
$ws.Range("C3").Value = 13.19193598703971
$ws.Range("D3").Value = 6.37239716649146
$ws.Range("F3").Value = 47.43714852535241
$ws.Range("G3").Value = 3.722575681919796
$ws.Range("J3").Value = 10.48446538920598
$ws.Range("M3").Value = 20.71457409415396
$ws.Range("N3").Value = 20.77557490723464

$ws.Range("B4").Value = 20.37345933097444
$ws.Range("C4").Value = 12.9469739269663
$ws.Range("D4").Value = 6.372625881170503
$ws.Range("F4").Value = 47.30249703531399
$ws.Range("G4").Value = 3.725447025913629
$ws.Range("J4").Value = 10.49949983060182
$ws.Range("M4").Value = 20.65381750008139
$ws.Range("N4").Value = 20.81186809679438

$ws.Range("B5").Value = 20.26696242165653
$ws.Range("C5").Value = 12.84698264906817
$ws.Range("D5").Value = 6.37276653022689
$ws.Range("F5").Value = 47.25079420846296
$ws.Range("G5").Value = 3.726652260587461
$ws.Range("J5").Value = 10.50592959237837
$ws.Range("M5").Value = 20.63035380065919
$ws.Range("N5").Value = 20.82715534648899

$ws.Range("B6").Value = 20.24932964446809
$ws.Range("C6").Value = 12.83037432381984
$ws.Range("D6").Value = 6.372792750549526
$ws.Range("F6").Value = 47.24240114956962
$ws.Range("G6").Value = 3.726854515436068
$ws.Range("J6").Value = 10.50701556099219
$ws.Range("M6").Value = 20.62653637404865
$ws.Range("N6").Value = 20.82972382264694

$ws.Range("B7").Value = 20.37201976545985
$ws.Range("C7").Value = 12.9456258505679
$ws.Range("D7").Value = 6.372627585904185
$ws.Range("F7").Value = 47.30178688705696
$ws.Range("G7").Value = 3.725463137660665
$ws.Range("J7").Value = 10.49958531713088
$ws.Range("M7").Value = 20.6534957944929
$ws.Range("N7").Value = 20.81207225198999

$ws.Range("B8").Value = 20.91969368572595
$ws.Range("C8").Value = 13.45226879281458
$ws.Range("D8").Value = 6.372315293947474
$ws.Range("F8").Value = 47.59158342489941
$ws.Range("G8").Value = 3.719631391759993
$ws.Range("J8").Value = 10.46945785640828
$ws.Range("M8").Value = 20.78378449237442
$ws.Range("N8").Value = 20.73854785243326

$ws.Range("B9").Value = 22.0055544572106
$ws.Range("C9").Value = 14.4278664995133
$ws.Range("D9").Value = 6.373200971469561
$ws.Range("F9").Value = 48.25738697385291
$ws.Range("G9").Value = 3.709292916147972
$ws.Range("J9").Value = 10.41993021436337
$ws.Range("M9").Value = 21.07871432633248
$ws.Range("N9").Value = 20.6100825669903

$ws.Range("B10").Value = 22.80137839644968
$ws.Range("C10").Value = 15.12487285325118
$ws.Range("D10").Value = 6.374767230511893
$ws.Range("F10").Value = 48.80363523764576
$ws.Range("G10").Value = 3.70235705488043
$ws.Range("J10").Value = 10.38935658845031
$ws.Range("M10").Value = 21.31815424806236
$ws.Range("N10").Value = 20.52529414129167

$ws.Range("B11").Value = 23.16115955341405
$ws.Range("C11").Value = 15.43602779064885
$ws.Range("D11").Value = 6.375679225272406
$ws.Range("F11").Value = 49.06393798038348
$ws.Range("G11").Value = 3.699343035067681
$ws.Range("J11").Value = 10.37670857289558
$ws.Range("M11").Value = 21.431715080178
$ws.Range("N11").Value = 20.48880894656754

$ws.Range("B12").Value = 23.29693363299974
$ws.Range("C12").Value = 15.55288260427604
$ws.Range("D12").Value = 6.376053311651118
$ws.Range("F12").Value = 49.16415007592496
$ws.Range("G12").Value = 3.698221845748693
$ws.Range("J12").Value = 10.37210016584496
$ws.Range("M12").Value = 21.47535650627091
$ws.Range("N12").Value = 20.47529332120874

$ws.Range("B13").Value = 23.26771525191903
$ws.Range("C13").Value = 15.52776091414728
$ws.Range("D13").Value = 6.375971466713192
$ws.Range("F13").Value = 49.14249557053434
$ws.Range("G13").Value = 3.698462419790185
$ws.Range("J13").Value = 10.37308461450129
$ws.Range("M13").Value = 21.46592961914301
$ws.Range("N13").Value = 20.4781907763476

$ws.Range("B14").Value = 23.17233993217854
$ws.Range("C14").Value = 15.44566166041468
$ws.Range("D14").Value = 6.375709425425172
$ws.Range("F14").Value = 49.07214994090513
$ws.Range("G14").Value = 3.699250391026279
$ws.Range("J14").Value = 10.37632580695023
$ws.Range("M14").Value = 21.4352928600685
$ws.Range("N14").Value = 20.48769098115277

$ws.Range("B15").Value = 23.1138547934038
$ws.Range("C15").Value = 15.39524338686017
$ws.Range("D15").Value = 6.375552661048217
$ws.Range("F15").Value = 49.02927316764326
$ws.Range("G15").Value = 3.699735666701672
$ws.Range("J15").Value = 10.3783347161082
$ws.Range("M15").Value = 21.41660923021754
$ws.Range("N15").Value = 20.49354928788075

$ws.Range("B16").Value = 22.77780816090297
$ws.Range("C16").Value = 15.10440839188114
$ws.Range("D16").Value = 6.374711647704717
$ws.Range("F16").Value = 48.78685669542914
$ws.Range("G16").Value = 3.702556856172624
$ws.Range("J16").Value = 10.39020851550389
$ws.Range("M16").Value = 21.31082370116311
$ws.Range("N16").Value = 20.52772055320145

$ws.Range("B17").Value = 22.5709744147259
$ws.Range("C17").Value = 14.92438461237293
$ws.Range("D17").Value = 6.374246839985251
$ws.Range("F17").Value = 48.64112758259761
$ws.Range("G17").Value = 3.704323613074234
$ws.Range("J17").Value = 10.39781538889093
$ws.Range("M17").Value = 21.24709576720583
$ws.Range("N17").Value = 20.54921809878249

$ws.Range("B18").Value = 22.4518061368499
$ws.Range("C18").Value = 14.82028932622915
$ws.Range("D18").Value = 6.373998269176023
$ws.Range("F18").Value = 48.55842349511301
$ws.Range("G18").Value = 3.705353097931055
$ws.Range("J18").Value = 10.40230927045066
$ws.Range("M18").Value = 21.21087952997042
$ws.Range("N18").Value = 20.56177917349828

$ws.Range("B19").Value = 22.4114274849361
$ws.Range("C19").Value = 14.78495383295922
$ws.Range("D19").Value = 6.373917330569745
$ws.Range("F19").Value = 48.53061458143421
$ws.Range("G19").Value = 3.705703951198495
$ws.Range("J19").Value = 10.40385119510524
$ws.Range("M19").Value = 21.19869347813209
$ws.Range("N19").Value = 20.56606582585495

$ws.Range("B20").Value = 22.59301430146203
$ws.Range("C20").Value = 14.94360633288875
$ws.Range("D20").Value = 6.374294376174232
$ws.Range("F20").Value = 48.65652565595654
$ws.Range("G20").Value = 3.704134164030714
$ws.Range("J20").Value = 10.39699334919344
$ws.Range("M20").Value = 21.25383454664555
$ws.Range("N20").Value = 20.54690933224578

$ws.Range("B21").Value = 23.2003677683303
$ws.Range("C21").Value = 15.46980351762925
$ws.Range("D21").Value = 6.375785613077039
$ws.Range("F21").Value = 49.09276809327263
$ws.Range("G21").Value = 3.699018398907097
$ws.Range("J21").Value = 10.37536887515959
$ws.Range("M21").Value = 21.44427452660083
$ws.Range("N21").Value = 20.48489237922781

$ws.Range("B22").Value = 23.59452375956836
$ws.Range("C22").Value = 15.80798506964681
$ws.Range("D22").Value = 6.376927731132737
$ws.Range("F22").Value = 49.3874143258556
$ws.Range("G22").Value = 3.695792361666695
$ws.Range("J22").Value = 10.36229171457792
$ws.Range("M22").Value = 21.57244509231595
$ws.Range("N22").Value = 20.44611273789007

$ws.Range("B23").Value = 23.38445572764288
$ws.Range("C23").Value = 15.62805178804612
$ws.Range("D23").Value = 6.376302817149322
$ws.Range("F23").Value = 49.22930332734401
$ws.Range("G23").Value = 3.697503461851132
$ws.Range("J23").Value = 10.36917467083943
$ws.Range("M23").Value = 21.50370860436618
$ws.Range("N23").Value = 20.46664963039101

$ws.Range("B24").Value = 22.58305085748264
$ws.Range("C24").Value = 14.9349180438597
$ws.Range("D24").Value = 6.374272826965559
$ws.Range("F24").Value = 48.6495608264426
$ws.Range("G24").Value = 3.704219771095453
$ws.Range("J24").Value = 10.39736461772803
$ws.Range("M24").Value = 21.25078662849146
$ws.Range("N24").Value = 20.54795249692886

$ws.Range("B25").Value = 21.71149886662935
$ws.Range("C25").Value = 14.16681686988157
$ws.Range("D25").Value = 6.372800778475423
$ws.Range("F25").Value = 48.06705616801953
$ws.Range("G25").Value = 3.7119732072414
$ws.Range("J25").Value = 10.43230705676025
$ws.Range("M25").Value = 20.99483223676532
$ws.Range("N25").Value = 20.64315202352579
